$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) The fixed smoke-test start date in B5 is replaced by a relative-date
#    placeholder token that the new python-appium-client test runner expands
#    at run time (today + 32 days).
$ws.Range("B5").Value = "<TODAY +32,+0,+0,'%m/%d/%Y'>"

# 2) Shrink the second illustration ("Grafik 3") so its right edge sits
#    further left - do this before widening column B so the anchor's
#    column/offset is computed against the original column geometry.
$ws.Shapes.Item(2).Width = 493.805

# 3) Column B needs to be wider now that it holds the longer placeholder text.
$ws.Columns.Item(2).ColumnWidth = 27.25

# 4) Restore the last-used cell selection on the sheet.
$ws.Range("E10").Select()
